$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 2.35
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 1.73
$ws.Range("N3").Value = 5
$ws.Range("W3").Value = 4.75
$ws.Range("AC3").Value = 4.5
$ws.Range("AM3").Value = 67
$ws.Range("AT3").Value = 1.91
$ws.Range("BA3").Value = 201

# Row 6
$ws.Range("G6").Value = 1.73
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = 2.38
$ws.Range("L6").Value = 5.5
$ws.Range("N6").Value = 9
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.7
$ws.Range("AA6").Value = 15
$ws.Range("AG6").Value = 800
$ws.Range("AJ6").Value = 17
$ws.Range("AM6").Value = 51
$ws.Range("AO6").Value = 9
$ws.Range("AW6").Value = 6.5
$ws.Range("AX6").Value = 29
$ws.Range("AY6").Value = 41
$ws.Range("BA6").Value = 151

# Row 7
$ws.Range("K7").Value = 2.62
$ws.Range("L7").Value = 6.1
$ws.Range("T7").Value = 3.6
$ws.Range("W7").Value = 9.5
$ws.Range("Z7").Value = 9.75
$ws.Range("AE7").Value = 16.5
$ws.Range("AH7").Value = 26
$ws.Range("AM7").Value = 50
$ws.Range("AN7").Value = 3.5
$ws.Range("AP7").Value = 13
$ws.Range("AS7").Value = 120
$ws.Range("AT7").Value = 3.6
$ws.Range("AU7").Value = 7.5
$ws.Range("AV7").Value = 50
$ws.Range("BB7").Value = 300

# Row 8
$ws.Range("K8").Value = 2.1
$ws.Range("L8").Value = 3.45
$ws.Range("O8").Value = 1.23
$ws.Range("T8").Value = 2.75
$ws.Range("AB8").Value = 23
$ws.Range("AH8").Value = 11
$ws.Range("AI8").Value = 17
$ws.Range("AM8").Value = 26
$ws.Range("AT8").Value = 2.75
$ws.Range("AU8").Value = 6.7
$ws.Range("AV8").Value = 55
$ws.Range("AX8").Value = 16
$ws.Range("AY8").Value = 22
$ws.Range("BA8").Value = 100
$ws.Range("BB8").Value = 250

# Row 12
$ws.Range("R12").Value = 3.35
$ws.Range("T12").Value = 4.25
$ws.Range("U12").Value = 1.78
$ws.Range("V12").Value = 1.93
$ws.Range("Z12").Value = 8.25
$ws.Range("AB12").Value = 23
$ws.Range("AD12").Value = 15.5
$ws.Range("AH12").Value = 50
$ws.Range("AK12").Value = 500
$ws.Range("AN12").Value = 3.5
$ws.Range("AO12").Value = 4.8
$ws.Range("AQ12").Value = 9.5
$ws.Range("AT12").Value = 4.25
$ws.Range("AW12").Value = 13.5
